# Generate Report for Handoff
#
# This updates the "Latest Handoff/Generate" timestamps for the row that
# corresponds to source file 9a19f432-d8c2-4a0c-90d9-f27498e5a3d1.md after a
# new handoff xliff was generated for it:
#   - Overview!G7  (Latest HO Xliff Generate Date) -> 2016-08-29 10:41:51
#   - zh-cn!H7     (Latest Handoff Datetime)        -> 2016-08-29 10:41:46
#   - de-de!H7     (Latest Handoff Datetime)        -> 2016-08-29 10:41:51

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G7").Value = "2016-08-29 10:41:51"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H7").Value = "2016-08-29 10:41:46"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H7").Value = "2016-08-29 10:41:51"
